$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2426.5454
$ws.Range("I19").Value = 2959.4
$ws.Range("J19").Value = 1982.5
$ws.Range("K19").Value = 2959.4
$ws.Range("L19").Value = 1982.5
$ws.Range("M19").Value = -2784.4
$ws.Range("N19").Value = -2332.5

$ws.Range("H106").Value = 2521.077
$ws.Range("I106").Value = 2442
$ws.Range("K106").Value = 2442
$ws.Range("M106").Value = -1811

$ws.Range("H113").Value = 83336830
$ws.Range("I113").Value = 25003750
$ws.Range("K113").Value = 25003750
$ws.Range("M113").Value = -25000496

$ws.Range("H131").Value = 5476.7
$ws.Range("I131").Value = 4635.6
$ws.Range("K131").Value = 13906.8
$ws.Range("M131").Value = -8866.800000000001

$ws.Range("H141").Value = 2500
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 2500
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 7500
$ws.Range("N141").Value = -17860
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1482.5
$ws.Range("I88").Value = 1239.8334
$ws.Range("J88").Value = 1664.5
$ws.Range("K88").Value = 1239.8334
$ws.Range("L88").Value = 1664.5
$ws.Range("M88").Value = -833.8334
$ws.Range("N88").Value = -2476.5

$ws.Range("H91").Value = 1482.5
$ws.Range("I91").Value = 1239.8334
$ws.Range("J91").Value = 1664.5
$ws.Range("K91").Value = 1239.8334
$ws.Range("L91").Value = 1664.5
$ws.Range("M91").Value = 164.1666
$ws.Range("N91").Value = -4472.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2545.6667
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 1818.5
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 1818.5
$ws.Range("M64").Value = -3775
$ws.Range("N64").Value = -2268.5

$ws.Range("H67").Value = 2545.6667
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 1818.5
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 1818.5
$ws.Range("M67").Value = -3220
$ws.Range("N67").Value = -3378.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1244.4333
$ws.Range("I7").Value = 109.71429
$ws.Range("J7").Value = 2237.3125
$ws.Range("K7").Value = 109.71429
$ws.Range("L7").Value = 2237.3125
$ws.Range("M7").Value = 3.285709999999995
$ws.Range("N7").Value = -2463.3125

$ws.Range("H107").Value = 679.5
$ws.Range("I107").Value = 593.38464
$ws.Range("K107").Value = 593.38464
$ws.Range("M107").Value = 1326.61536

$ws.Range("H122").Value = 4719.294
$ws.Range("I122").Value = 2578
$ws.Range("J122").Value = 5178.143
$ws.Range("K122").Value = 7734
$ws.Range("L122").Value = 15534.429
$ws.Range("M122").Value = -5284
$ws.Range("N122").Value = -20434.429

$ws.Range("H141").Value = 96533.125
$ws.Range("J141").Value = 103799.91
$ws.Range("L141").Value = 103799.91
$ws.Range("N141").Value = -114159.91

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1942.5555
$ws.Range("I34").Value = 489.8
$ws.Range("J34").Value = 2501.3076
$ws.Range("K34").Value = 1469.4
$ws.Range("L34").Value = 7503.9228
$ws.Range("M34").Value = -1385.4
$ws.Range("N34").Value = -7671.9228

$ws.Range("H113").Value = 2502.3333
$ws.Range("J113").Value = 3134.75
$ws.Range("L113").Value = 9404.25
$ws.Range("N113").Value = -13744.25

$ws.Range("H131").Value = 1504.9231
$ws.Range("I131").Value = 1351.1666
$ws.Range("J131").Value = 3350
$ws.Range("K131").Value = 4053.4998
$ws.Range("L131").Value = 10050
$ws.Range("M131").Value = 986.5001999999999
$ws.Range("N131").Value = -20130

$ws.Range("H132").Value = 2317.6667
$ws.Range("J132").Value = 2106.3333
$ws.Range("L132").Value = 18956.9997
$ws.Range("N132").Value = -24016.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3453.0278
$ws.Range("I102").Value = 2264.84
$ws.Range("K102").Value = 2264.84
$ws.Range("M102").Value = -642.8400000000001

$ws.Range("H123").Value = 60225.332
$ws.Range("J123").Value = 60225.332
$ws.Range("L123").Value = 60225.332
$ws.Range("N123").Value = -65125.332

$ws.Range("H124").Value = 105026.164
$ws.Range("J124").Value = 105026.164
$ws.Range("L124").Value = 105026.164
$ws.Range("N124").Value = -114846.164

$ws.Range("H126").Value = 6142.857
$ws.Range("I126").Value = 7666.6665
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 22999.9995
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -20529.9995
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 75940.36
$ws.Range("I7").Value = 3808
$ws.Range("J7").Value = 130039.625
$ws.Range("K7").Value = 3808
$ws.Range("L7").Value = 130039.625
$ws.Range("M7").Value = -3696
$ws.Range("N7").Value = -130263.625

$ws.Range("H40").Value = 4233.636
$ws.Range("I40").Value = 3231.5715
$ws.Range("J40").Value = 5987.25
$ws.Range("K40").Value = 3231.5715
$ws.Range("L40").Value = 5987.25
$ws.Range("M40").Value = -3095.5715
$ws.Range("N40").Value = -6259.25

$ws.Range("H68").Value = 977.1818
$ws.Range("I68").Value = 546.8570999999999
$ws.Range("K68").Value = 546.8570999999999
$ws.Range("M68").Value = 202.1429000000001

$ws.Range("H71").Value = 977.1818
$ws.Range("I71").Value = 546.8570999999999
$ws.Range("K71").Value = 2734.2855
$ws.Range("M71").Value = 1009.7145

$ws.Range("H122").Value = 5414.6665
$ws.Range("I122").Value = 4664.1113
$ws.Range("K122").Value = 13992.3339
$ws.Range("M122").Value = -11542.3339

$ws.Range("H126").Value = 75940.36
$ws.Range("I126").Value = 3808
$ws.Range("J126").Value = 130039.625
$ws.Range("K126").Value = 11424
$ws.Range("L126").Value = 390118.875
$ws.Range("M126").Value = -8954
$ws.Range("N126").Value = -395058.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5665.3
$ws.Range("I126").Value = 5715.125
$ws.Range("J126").Value = 5466
$ws.Range("K126").Value = 17145.375
$ws.Range("L126").Value = 16398
$ws.Range("M126").Value = -14675.375
$ws.Range("N126").Value = -21338

$ws.Range("H133").Value = 98850
$ws.Range("J133").Value = 98850
$ws.Range("L133").Value = 98850
$ws.Range("N133").Value = -108970
